# Update Excel COM-interop script implementing the commit
# "Update countries & provincias Spain"
#
# This updates the COVID-19 country statistics table on worksheet "Pais":
#  - refreshes the "last updated" timestamp string
#  - refreshes several countries' case/recovery/death counters
#  - re-sorts the Malta / Santo Tome y Principe / Siria block so that
#    Siria (whose case count grew past Malta's) now appears first

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header: "last updated" timestamp -------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 4 de Agosto de 2020 a las 21:49"

# --- Row 4: Estados Unidos --------------------------------------------------
$ws.Cells.Item(4, 2).Value = 4895491
$ws.Cells.Item(4, 3).Value = 31575
$ws.Cells.Item(4, 4).Value = 2462542
$ws.Cells.Item(4, 5).Value = 2273176
$ws.Cells.Item(4, 7).Value = 845
$ws.Cells.Item(4, 8).Value = 159773

# --- Row 6: India ------------------------------------------------------------
$ws.Cells.Item(6, 2).Value = 1906613
$ws.Cells.Item(6, 3).Value = 51282
$ws.Cells.Item(6, 4).Value = 1281660
$ws.Cells.Item(6, 5).Value = 585133

# --- Row 23: Francia ----------------------------------------------------------
$ws.Cells.Item(23, 5).Value = 79872
$ws.Cells.Item(23, 7).Value = 2
$ws.Cells.Item(23, 8).Value = 30296

# --- Row 52: Barein ------------------------------------------------------------
$ws.Cells.Item(52, 5).Value = 2677
$ws.Cells.Item(52, 7).Value = 1
$ws.Cells.Item(52, 8).Value = 151

# --- Row 57: Afganistan --------------------------------------------------------
$ws.Cells.Item(57, 2).Value = 36782
$ws.Cells.Item(57, 3).Value = 35
$ws.Cells.Item(57, 5).Value = 9825

# --- Row 62: Uzbekistan --------------------------------------------------------
$ws.Cells.Item(62, 2).Value = 27047
$ws.Cells.Item(62, 3).Value = 981
$ws.Cells.Item(62, 4).Value = 18051
$ws.Cells.Item(62, 5).Value = 8831
$ws.Cells.Item(62, 7).Value = 8
$ws.Cells.Item(62, 8).Value = 165

# --- Row 89: Guayana Francesa ---------------------------------------------------
$ws.Cells.Item(89, 2).Value = 7998
$ws.Cells.Item(89, 3).Value = 50
$ws.Cells.Item(89, 4).Value = 6873
$ws.Cells.Item(89, 5).Value = 1080
$ws.Cells.Item(89, 7).Value = 1
$ws.Cells.Item(89, 8).Value = 45

# --- Row 108: Malaui -------------------------------------------------------------
$ws.Cells.Item(108, 2).Value = 4361
$ws.Cells.Item(108, 3).Value = 89
$ws.Cells.Item(108, 4).Value = 2047
$ws.Cells.Item(108, 5).Value = 2186
$ws.Cells.Item(108, 7).Value = 5
$ws.Cells.Item(108, 8).Value = 128

# --- Rows 153-155: Siria moves above Malta / Santo Tome y Principe ----------------
# Row 153 becomes Siria with its updated counters
$ws.Cells.Item(153, 1).Value = "Siria"
$ws.Cells.Item(153, 2).Value = 892
$ws.Cells.Item(153, 3).Value = 45
$ws.Cells.Item(153, 4).Value = 283
$ws.Cells.Item(153, 5).Value = 563
$ws.Cells.Item(153, 6).Value = 0
$ws.Cells.Item(153, 7).Value = 0
$ws.Cells.Item(153, 8).Value = 46

# Row 154 becomes Malta (previous row-153 data)
$ws.Cells.Item(154, 1).Value = "Malta"
$ws.Cells.Item(154, 2).Value = 890
$ws.Cells.Item(154, 3).Value = 16
$ws.Cells.Item(154, 4).Value = 666
$ws.Cells.Item(154, 5).Value = 215
$ws.Cells.Item(154, 6).Value = 0
$ws.Cells.Item(154, 7).Value = 0
$ws.Cells.Item(154, 8).Value = 9

# Row 155 becomes Santo Tome y Principe (previous row-154 data)
$ws.Cells.Item(155, 1).Value = "Santo Tome y Principe"
$ws.Cells.Item(155, 2).Value = 874
$ws.Cells.Item(155, 3).Value = 0
$ws.Cells.Item(155, 4).Value = 787
$ws.Cells.Item(155, 5).Value = 72
$ws.Cells.Item(155, 6).Value = 0
$ws.Cells.Item(155, 7).Value = 0
$ws.Cells.Item(155, 8).Value = 15
